$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J3").Value = 1.11
$ws.Range("K3").Value = 6.5
$ws.Range("N3").Value = 2.6
$ws.Range("O3").Value = 1.48
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.38
$ws.Range("L5").Value = 1.53
$ws.Range("M5").Value = 2.38
$ws.Range("N5").Value = 2.7
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 1.57
$ws.Range("AE6").Value = 9.5
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 19
$ws.Range("I6").Value = 2.05
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.63
$ws.Range("X6").Value = 29
$ws.Range("AB8").Value = 28
$ws.Range("AD8").Value = 5.1
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 14
$ws.Range("AH8").Value = 50
$ws.Range("AI8").Value = 100
$ws.Range("H8").Value = 2.32
$ws.Range("J8").Value = 1.26
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 1.72
$ws.Range("N8").Value = 3.95
$ws.Range("O8").Value = 1.21
$ws.Range("P8").Value = 1.88
$ws.Range("Q8").Value = 1.82
$ws.Range("R8").Value = 2.77
$ws.Range("S8").Value = 1.39
$ws.Range("T8").Value = 5.3
$ws.Range("U8").Value = 13.5
$ws.Range("V8").Value = 14
$ws.Range("AI9").Value = 51
$ws.Range("AJ9").Value = 401
$ws.Range("G9").Value = 1.65
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 5.25
$ws.Range("J9").Value = 1.06
$ws.Range("K9").Value = 10
$ws.Range("N9").Value = 2.05
$ws.Range("O9").Value = 1.75
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.63
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 1.73
$ws.Range("X9").Value = 15
$ws.Range("Z9").Value = 9
$ws.Range("AA15").Value = 9.5
$ws.Range("AC15").Value = 29
$ws.Range("AD15").Value = 19
$ws.Range("AE15").Value = 26
$ws.Range("AH15").Value = 26
$ws.Range("AJ15").Value = 81
$ws.Range("G15").Value = 1.75
$ws.Range("H15").Value = 4.5
$ws.Range("I15").Value = 3.8
$ws.Range("J15").Value = 1.01
$ws.Range("K15").Value = 26
$ws.Range("P15").Value = 1.19
$ws.Range("T15").Value = 13
$ws.Range("Z15").Value = 26
$ws.Range("AJ16").Value = 151
$ws.Range("N16").Value = 1.57
$ws.Range("P16").Value = 1.25
$ws.Range("Q16").Value = 3.5
$ws.Range("AG17").Value = 67
$ws.Range("N17").Value = 1.47
$ws.Range("P17").Value = 1.22
$ws.Range("Z17").Value = 17
$ws.Range("AD19").Value = 10
$ws.Range("AE19").Value = 19
$ws.Range("AH19").Value = 34
$ws.Range("AJ19").Value = 351
$ws.Range("G19").Value = 2.05
$ws.Range("I19").Value = 3.7
$ws.Range("R19").Value = 1.91
$ws.Range("S19").Value = 1.8
$ws.Range("T19").Value = 6.5
$ws.Range("U19").Value = 9.5
$ws.Range("W19").Value = 19
$ws.Range("AA22").Value = 6
$ws.Range("AD22").Value = 9.5
$ws.Range("AE22").Value = 15
$ws.Range("AF22").Value = 11
$ws.Range("AG22").Value = 29
$ws.Range("AH22").Value = 23
$ws.Range("AI22").Value = 29
$ws.Range("G22").Value = 2.4
$ws.Range("H22").Value = 3.1
$ws.Range("I22").Value = 2.7
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 1.85
$ws.Range("T22").Value = 9
$ws.Range("U22").Value = 12
$ws.Range("V22").Value = 10
$ws.Range("W22").Value = 23
$ws.Range("X22").Value = 21
$ws.Range("Y22").Value = 29
$ws.Range("AG25").Value = 81
$ws.Range("G25").Value = 1.25
$ws.Range("H25").Value = 6.5
$ws.Range("I25").Value = 7
$ws.Range("L25").Value = 1.1
$ws.Range("M25").Value = 6.5
$ws.Range("N25").Value = 1.33
$ws.Range("O25").Value = 3.25
$ws.Range("R25").Value = 1.67
$ws.Range("S25").Value = 2.1
$ws.Range("T25").Value = 12
$ws.Range("V25").Value = 9.5
$ws.Range("W25").Value = 9.5
$ws.Range("AB29").Value = 16.5
$ws.Range("H29").Value = 3.4
$ws.Range("I29").Value = 4.65
$ws.Range("L29").Value = 1.32
$ws.Range("M29").Value = 2.82
$ws.Range("X29").Value = 15
